# Commit: "Added create new table" - clears the sample/demo data that was
# filled in for rows 2-8 so users can build their own schedule from a
# blank template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @(
    "A2","B2","C2","E2","F2","G2","O2",
    "A3","E3","F3","G3","H3",
    "A4","B4","C4","D4","E4","F4","G4","H4","O4",
    "A5","B5","D5","F5",
    "A6","B6","C6","D6","E6","F6","G6","H6",
    "A7","B7","C7","E7","F7","G7","H7",
    "A8","C8","D8","E8","F8","G8","H8"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = $null
}
